# DC-Colos.xlsx update: insert a new "Agra, India" (AGR) colo row just
# before the existing "Ashburn, VA" (IAD) row, pushing IAD..YHZ (old rows
# 274-332) down by one row (new rows 275-333).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 274; everything from the old row 274
# downward (IAD, ATL, BOS, ... YHZ) shifts down to row 275.. onward, and
# the sheet dimension grows from F332 to F333 automatically.
$ws.Rows(274).Insert()

# The freshly inserted row lost the header-style border that column A
# carries throughout the table (Insert() only carries the font/alignment
# forward). Copy the formatting from the row above (still "CTU", the row
# that used to be 273 and is unaffected by the insert) onto the new row
# so A274 again matches the rest of the colo-code column.
$ws.Range("A273").Copy()
$ws.Range("A274").PasteSpecial(-4122)

# Populate the new row with the Agra, India colo data.
$ws.Range("A274").Value = "AGR"
$ws.Range("B274").Value = "Agra, India"
$ws.Range("C274").Value = "Asia"
$ws.Range("D274").Value = "Agra"
$ws.Range("E274").Value = "India"
$ws.Range("F274").Value = "IN"
